$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.899.42'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.75%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.514.60'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.71%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.02'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.45%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.08'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.88%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.513.23'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.78%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.75%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.143'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.94%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.96'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +5.56%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.81%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.110.71'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.96%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '31.94'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.35%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.512.80'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.79%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.989.60'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.56%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.76'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +8.76%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.46%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.37'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '438.43'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.76%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.65%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.46'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.654.02'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.76%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.78'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.55%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.31'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.96%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.74%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.62'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.39%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.60%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.10%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.21%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.97'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.50%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.13%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.04'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.63%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.07%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '174.10'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.64%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0896'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.09%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.16%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -10.03%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.895'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.45%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '46.18'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.67%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '27.84'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -7.79%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.28'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.24%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.56%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.47'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.74%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.994'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.86%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.249'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.02%  '
